$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auditoria / revision de no conformidades:
# Todas las no conformidades quedan marcadas como "Cerrada" y se
# registra la fecha real de cierre (columna E). La fila 6 no aplica
# fecha, por lo que se indica "N/a".

$ws.Range("F4:F9").Value = "Cerrada"

$ws.Range("E4").Value = $ws.Range("D4").Value2
$ws.Range("E5").Value = $ws.Range("D5").Value2
$ws.Range("E6").Value = "N/a"
$ws.Range("E7").Value = $ws.Range("D7").Value2
$ws.Range("E8").Value = $ws.Range("D8").Value2
$ws.Range("E9").Value = $ws.Range("D9").Value2

[void]$ws.Range("D9").Select()
